# Insert a new data row at row 495 (pushes existing rows 495:593 down to 496:594)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(495).Insert()

# Populate the newly inserted row with its values
$ws.Range("A495").Value = 6
$ws.Range("B495").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C495").Value = "Metropolitana"
$ws.Range("D495").Value = 45015
$ws.Range("E495").Value = 13
$ws.Range("F495").Value = 100112043
$ws.Range("G495").Value = "Pepino ensalada"
$ws.Range("H495").Value = "Sin especificar"
$ws.Range("I495").Value = "Primera"
$ws.Range("J495").Value = 580
$ws.Range("K495").Value = 5000
$ws.Range("L495").Value = 6000
$ws.Range("M495").Value = 5448
$ws.Range("N495").Value = "`$/caja 60 unidades"
$ws.Range("O495").Value = "Región de Arica y Parinacota"
$ws.Range("P495").Value = 91
$ws.Range("Q495").Value = 60
$ws.Range("R495").Value = "Hortaliza"
